$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = 1398800851.3300049
$ws.Range("E12").Value = 1361974149

# Row 13
$ws.Range("D13").Value = 338965688.90999979
$ws.Range("E13").Value = 332114255.60000002

# Row 14
$ws.Range("D14").Value = -45752811.059999987
$ws.Range("E14").Value = 637985.22

# Row 16
$ws.Range("D16").Value = -50601311.959999993

# Row 18 - D18 becomes a formula (was a plain value)
$ws.Range("D18").Formula = "=SUM(D12:D17)"

# Row 19
$ws.Range("D19").Value = -383099999.99999988
$ws.Range("E19").Value = -384700000

# Row 21 - D21 becomes a formula (was a plain value)
$ws.Range("D21").Formula = "=SUM(D18:D20)"

# Row 22
$ws.Range("D22").Value = -32201025

# Row 26
$ws.Range("D26").Value = 1009991810.1331247
$ws.Range("E26").Value = 101863404

Write-Host "Edit applied"
